$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (ID + Pesquisa) below the existing table
$ws.Range("A4").Value = "G_00003"
$ws.Range("B4").Value = "RSI"

# Match the formatting (thin border) used by the previous data row
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-assert the values in case PasteSpecial touched them
$ws.Range("A4").Value = "G_00003"
$ws.Range("B4").Value = "RSI"

$ws.Range("B16").Select()
